# issue #5: stock data output to json file
#
# The "股票" (stock) sheet gets a new "property_category" column inserted
# right before the existing "date" column, with the literal value "stock"
# for the (single) data row. This pushes date / legislator_name /
# legislator_id one column to the right (H->I, I->J, J->K).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column at H (shifting the old H:J -> I:K, xlShiftToRight = -4161)
$ws.Range("H1:H2").Insert(-4161)

# Populate the newly inserted column.
$ws.Range("H1").Value = "property_category"
$ws.Range("H2").Value = "stock"
